# Applies the corrected IFRS figures for 2014-2018 (rows 2-6) and clears the
# stale 2019E-2021E forecast rows (7-9), matching the upstream data fix.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$rowValues = @{
    "D2" = 331
    "E2" = -54
    "F2" = -54
    "G2" = -44
    "H2" = -44
    "I2" = -39
    "J2" = -5
    "K2" = 848
    "L2" = 227
    "M2" = 621
    "N2" = 623
    "O2" = -2
    "P2" = 774
    "Q2" = -69
    "R2" = -150
    "S2" = 133
    "T2" = 20
    "U2" = -89
    "V2" = 141
    "W2" = -16.17
    "X2" = -13.28
    "Y2" = -6.56
    "Z2" = -5.39
    "AA2" = 36.57
    "AB2" = -21.76
    "AC2" = -42
    "AD2" = -16.16
    "AE2" = 612
    "AF2" = 1.11
    "AG2" = 0
    "AH2" = 0
    "AI2" = 0
    "AJ2" = 101839304
}
foreach ($ref in $rowValues.Keys) {
    $ws.Range($ref).Value = $rowValues[$ref]
}

# Row 3
$rowValues = @{
    "D3" = 403
    "E3" = -9
    "F3" = -9
    "G3" = 12
    "H3" = 16
    "I3" = 28
    "J3" = -12
    "K3" = 852
    "L3" = 175
    "M3" = 677
    "N3" = 692
    "O3" = -15
    "P3" = 403
    "Q3" = 18
    "R3" = 101
    "S3" = -33
    "T3" = 54
    "U3" = -35
    "V3" = 87
    "W3" = -2.15
    "X3" = 3.85
    "Y3" = 4.2
    "Z3" = 1.83
    "AA3" = 25.8
    "AB3" = 68.09
    "AC3" = 26
    "AD3" = 100.05
    "AE3" = 653
    "AF3" = 4.05
    "AG3" = 0
    "AH3" = 0
    "AI3" = 0
    "AJ3" = 106028154
}
foreach ($ref in $rowValues.Keys) {
    $ws.Range($ref).Value = $rowValues[$ref]
}

# Row 4
$rowValues = @{
    "D4" = 441
    "E4" = -2
    "F4" = -2
    "G4" = -41
    "H4" = -48
    "I4" = -39
    "J4" = -8
    "K4" = 1287
    "L4" = 506
    "M4" = 781
    "N4" = 804
    "O4" = -23
    "P4" = 404
    "Q4" = 4
    "R4" = -355
    "S4" = 434
    "T4" = 17
    "U4" = -12
    "V4" = 370
    "W4" = -0.5600000000000001
    "X4" = -10.8
    "Y4" = -5.24
    "Z4" = -4.46
    "AA4" = 64.73
    "AB4" = 95.77
    "AC4" = -37
    "AD4" = -87.75
    "AE4" = 756
    "AF4" = 4.28
    "AG4" = 0
    "AH4" = 0
    "AI4" = 0
    "AJ4" = 106425300
}
foreach ($ref in $rowValues.Keys) {
    $ws.Range($ref).Value = $rowValues[$ref]
}

# Row 5
$rowValues = @{
    "D5" = 461
    "E5" = -15
    "F5" = -15
    "G5" = -165
    "H5" = -169
    "I5" = -164
    "J5" = -5
    "K5" = 1572
    "L5" = 804
    "M5" = 768
    "N5" = 768
    "P5" = 423
    "Q5" = -89
    "R5" = -339
    "S5" = 405
    "T5" = 21
    "U5" = -110
    "V5" = 517
    "W5" = -3.16
    "X5" = -36.64
    "Y5" = -20.82
    "Z5" = -11.82
    "AA5" = 104.68
    "AB5" = 83.63
    "AC5" = -152
    "AD5" = -22.23
    "AE5" = 690
    "AF5" = 4.91
    "AG5" = 0
    "AH5" = 0
    "AI5" = 0
    "AJ5" = 111365254
}
foreach ($ref in $rowValues.Keys) {
    $ws.Range($ref).Value = $rowValues[$ref]
}

# Row 6
$rowValues = @{
    "D6" = 454
    "E6" = -8
    "F6" = -8
    "G6" = -10
    "H6" = -12
    "I6" = -12
    "K6" = 2171
    "L6" = 666
    "M6" = 1505
    "N6" = 1505
    "P6" = 539
    "Q6" = 13
    "R6" = -581
    "S6" = 546
    "T6" = 10
    "U6" = 3
    "V6" = 380
    "W6" = -1.82
    "X6" = -2.68
    "Y6" = -1.07
    "Z6" = -0.65
    "AA6" = 44.22
    "AB6" = 175.53
    "AC6" = -9
    "AD6" = -209.77
    "AE6" = 1060
    "AF6" = 1.86
    "AG6" = 0
    "AH6" = 0
    "AI6" = 0
    "AJ6" = 141964790
}
foreach ($ref in $rowValues.Keys) {
    $ws.Range($ref).Value = $rowValues[$ref]
}

# O5 no longer reported for this year - drop it entirely
$ws.Range("O5").ClearContents()

# 2019E/2020E/2021E estimate rows were erroneous - clear all figures, keep labels
$ws.Range("D7:AJ9").ClearContents()
